{"js": "// 1) \"M\u00fc\u015fteri (Customer)\" bullet: insert \"(Cache)\" right after \"g\u00f6r\u00fcnt\u00fcleyebilir\"\n//    that precedes \", sipari\u015f olu\u015fturabilir\".\nconst custResults = context.document.body.search(\n  \"g\u00f6r\u00fcnt\u00fcleyebilir, sipari\u015f olu\u015fturabilir\",\n  { matchCase: true }\n);\ncustResults.load(\"text\");\nawait context.sync();\n\nif (custResults.items.length > 0) {\n  custResults.items[0].insertText(\n    \"g\u00f6r\u00fcnt\u00fcleyebilir(Cache), sipari\u015f olu\u015fturabilir\",\n    \"Replace\"\n  );\n}\n\n// 2) \"Admin\" bullet: insert \"(cache, sayfalama)\" right after \"listeleyebilir\"\n//    that precedes \", kullan\u0131c\u0131 ve restoran bilgilerini\".\nconst adminResults = context.document.body.search(\n  \"listeleyebilir, kullan\u0131c\u0131 ve restoran\",\n  { matchCase: true }\n);\nadminResults.load(\"text\");\nawait context.sync();\n\nif (adminResults.items.length > 0) {\n  adminResults.items[0].insertText(\n    \"listeleyebilir(cache, sayfalama), kullan\u0131c\u0131 ve restoran\",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n\n// 3) Add two new bullet items (\"Caching\" and \"Sayfalama - Pageable\") right\n//    after the \"Global Exception Handling (...)\" bullet, before the\n//    \"Swagger kullanarak API d\u00f6k\u00fcmantasyonu olu\u015fturma.\" bullet. Search for\n//    text unique to the end of that specific bullet (there are other\n//    \"Global Exception Handling\" mentions elsewhere in the document).\nconst techResults = context.document.body.search(\n  \"@ExceptionHandler)\",\n  { matchCase: true }\n);\ntechResults.load(\"text\");\nawait context.sync();\n\nif (techResults.items.length > 0) {\n  const techParagraph = techResults.items[0].paragraphs.getFirst();\n  const cachingParagraph = techParagraph.insertParagraph(\"Caching\", \"After\");\n  cachingParagraph.insertParagraph(\"Sayfalama - Pageable\", \"After\");\n}\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) \"M\u00fc\u015fteri (Customer)\" bullet: insert \"(Cache)\" right after \"g\u00f6r\u00fcnt\u00fcleyebilir\"\n#    that precedes \", sipari\u015f olu\u015fturabilir\".\n$rng1 = $d.Content\n$rng1.Find.Execute(\n    \"g\u00f6r\u00fcnt\u00fcleyebilir, sipari\u015f olu\u015fturabilir\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"g\u00f6r\u00fcnt\u00fcleyebilir(Cache), sipari\u015f olu\u015fturabilir\",\n    2\n) | Out-Null\n\n# 2) \"Admin\" bullet: insert \"(cache, sayfalama)\" right after \"listeleyebilir\"\n#    that precedes \", kullan\u0131c\u0131 ve restoran bilgilerini\".\n$rng2 = $d.Content\n$rng2.Find.Execute(\n    \"listeleyebilir, kullan\u0131c\u0131 ve restoran\",\n    $false, $false, $false, $false, $false, $true, 1, $false,\n    \"listeleyebilir(cache, sayfalama), kullan\u0131c\u0131 ve restoran\",\n    2\n) | Out-Null\n\n# 3) Add two new bullet items (\"Caching\" and \"Sayfalama - Pageable\") right\n#    after the \"Global Exception Handling (...)\" bullet, before the\n#    \"Swagger kullanarak API d\u00f6k\u00fcmantasyonu olu\u015fturma.\" bullet. Locate that\n#    specific bullet by text unique to its end (other, unrelated \"Global\n#    Exception Handling\" mentions exist elsewhere in the document).\n$targetIdx = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*@ExceptionHandler)*\") {\n        $targetIdx = $i\n        break\n    }\n}\n\nif ($targetIdx -gt 0) {\n    $d.Paragraphs($targetIdx).Range.InsertParagraphAfter()\n    $d.Paragraphs($targetIdx + 1).Range.Text = \"Caching\"\n\n    $d.Paragraphs($targetIdx + 1).Range.InsertParagraphAfter()\n    $d.Paragraphs($targetIdx + 2).Range.Text = \"Sayfalama - Pageable\"\n}\n"}
